# Auto-generated edit script: refreshes cached Universalis market-price-derived
# values (currentAveragePrice*, Leve/Profit columns) per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 11759.8
$ws.Range("I13").Value = 4950
$ws.Range("K13").Value = 4950
$ws.Range("M13").Value = -4781
$ws.Range("H19").Value = 4168122.5
$ws.Range("I19").Value = 1194.2222
$ws.Range("J19").Value = 6668279.5
$ws.Range("K19").Value = 1194.2222
$ws.Range("L19").Value = 6668279.5
$ws.Range("M19").Value = -1019.2222
$ws.Range("N19").Value = -6668629.5
$ws.Range("H28").Value = 5358.1763
$ws.Range("I28").Value = 949.5
$ws.Range("K28").Value = 949.5
$ws.Range("M28").Value = -464.5
$ws.Range("H33").Value = 243.28572
$ws.Range("I33").Value = 254.45
$ws.Range("K33").Value = 254.45
$ws.Range("M33").Value = -25.44999999999999
$ws.Range("H39").Value = 1190.1
$ws.Range("I39").Value = 56.42857
$ws.Range("J39").Value = 3835.3333
$ws.Range("K39").Value = 169.28571
$ws.Range("L39").Value = 11505.9999
$ws.Range("M39").Value = 126.71429
$ws.Range("N39").Value = -12097.9999
$ws.Range("H41").Value = 946.7857
$ws.Range("I41").Value = 1570.5714
$ws.Range("J41").Value = 323
$ws.Range("K41").Value = 1570.5714
$ws.Range("L41").Value = 323
$ws.Range("M41").Value = -1130.5714
$ws.Range("N41").Value = -1203
$ws.Range("H43").Value = 4350.75
$ws.Range("I43").Value = 4258
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 4258
$ws.Range("L43").Value = 5000
$ws.Range("M43").Value = -4189
$ws.Range("N43").Value = -5138
$ws.Range("H62").Value = 4636.421
$ws.Range("J62").Value = 10319.8
$ws.Range("L62").Value = 10319.8
$ws.Range("N62").Value = -11567.8
$ws.Range("H65").Value = 4636.421
$ws.Range("J65").Value = 10319.8
$ws.Range("L65").Value = 51599
$ws.Range("N65").Value = -57839
$ws.Range("H74").Value = 98569.5
$ws.Range("I74").Value = 116493.4
$ws.Range("K74").Value = 116493.4
$ws.Range("M74").Value = -115557.4
$ws.Range("H76").Value = 7847.727
$ws.Range("I76").Value = 10078.714
$ws.Range("J76").Value = 3943.5
$ws.Range("K76").Value = 10078.714
$ws.Range("L76").Value = 3943.5
$ws.Range("M76").Value = -9763.714
$ws.Range("N76").Value = -4573.5
$ws.Range("H77").Value = 98569.5
$ws.Range("I77").Value = 116493.4
$ws.Range("K77").Value = 582467
$ws.Range("M77").Value = -577787
$ws.Range("H79").Value = 7847.727
$ws.Range("I79").Value = 10078.714
$ws.Range("J79").Value = 3943.5
$ws.Range("K79").Value = 10078.714
$ws.Range("L79").Value = 3943.5
$ws.Range("M79").Value = -8986.714
$ws.Range("N79").Value = -6127.5
$ws.Range("H86").Value = 6798.846
$ws.Range("I86").Value = 3459.4443
$ws.Range("K86").Value = 3459.4443
$ws.Range("M86").Value = -2336.4443
$ws.Range("H89").Value = 6798.846
$ws.Range("I89").Value = 3459.4443
$ws.Range("K89").Value = 17297.2215
$ws.Range("M89").Value = -11681.2215
$ws.Range("H92").Value = 2740.5833
$ws.Range("I92").Value = 1608.25
$ws.Range("J92").Value = 5005.25
$ws.Range("K92").Value = 1608.25
$ws.Range("L92").Value = 5005.25
$ws.Range("M92").Value = -360.25
$ws.Range("N92").Value = -7501.25
$ws.Range("H98").Value = 3474171.2
$ws.Range("I98").Value = 3678122.5
$ws.Range("J98").Value = 7000
$ws.Range("K98").Value = 3678122.5
$ws.Range("L98").Value = 7000
$ws.Range("M98").Value = -3676624.5
$ws.Range("N98").Value = -9996
$ws.Range("H100").Value = 4670.353
$ws.Range("I100").Value = 2720
$ws.Range("K100").Value = 2720
$ws.Range("M100").Value = -2179
$ws.Range("H106").Value = 7831.5835
$ws.Range("I106").Value = 9888.111000000001
$ws.Range("K106").Value = 9888.111000000001
$ws.Range("M106").Value = -9257.111000000001
$ws.Range("H107").Value = 2799.125
$ws.Range("I107").Value = 599.75
$ws.Range("J107").Value = 4998.5
$ws.Range("K107").Value = 599.75
$ws.Range("L107").Value = 4998.5
$ws.Range("M107").Value = 1320.25
$ws.Range("N107").Value = -8838.5
$ws.Range("H109").Value = 71871.664
$ws.Range("J109").Value = 99997
$ws.Range("L109").Value = 99997
$ws.Range("N109").Value = -102771
$ws.Range("H122").Value = 3474171.2
$ws.Range("I122").Value = 3678122.5
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 11034367.5
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -11031917.5
$ws.Range("N122").Value = -25900
$ws.Range("H124").Value = 99998.5
$ws.Range("J124").Value = 99998.5
$ws.Range("L124").Value = 99998.5
$ws.Range("N124").Value = -109818.5
$ws.Range("H125").Value = 8046.2
$ws.Range("I125").Value = 6057.75
$ws.Range("K125").Value = 54519.75
$ws.Range("M125").Value = -52059.75
$ws.Range("H132").Value = 2397.0625
$ws.Range("I132").Value = 2173.3667
$ws.Range("K132").Value = 6520.1001
$ws.Range("M132").Value = -3990.1001
$ws.Range("H137").Value = 3553.8635
$ws.Range("I137").Value = 4273.8335
$ws.Range("K137").Value = 12821.5005
$ws.Range("M137").Value = -10271.5005
$ws.Range("H141").Value = 14713528
$ws.Range("I141").Value = 17860466
$ws.Range("K141").Value = 53581398
$ws.Range("M141").Value = -53576218

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1706.5625
$ws.Range("I2").Value = 916
$ws.Range("K2").Value = 916
$ws.Range("M2").Value = -803
$ws.Range("H61").Value = 7007578.5
$ws.Range("I61").Value = 7699909.5
$ws.Range("K61").Value = 7699909.5
$ws.Range("M61").Value = -7699697.5
$ws.Range("H116").Value = 1706.5625
$ws.Range("I116").Value = 916
$ws.Range("K116").Value = 916
$ws.Range("M116").Value = 1378
$ws.Range("H132").Value = 2947522
$ws.Range("I132").Value = 5097.7036
$ws.Range("J132").Value = 14296873
$ws.Range("K132").Value = 15293.1108
$ws.Range("L132").Value = 42890619
$ws.Range("M132").Value = -12763.1108
$ws.Range("N132").Value = -42895679
$ws.Range("H136").Value = 7007578.5
$ws.Range("I136").Value = 7699909.5
$ws.Range("K136").Value = 23099728.5
$ws.Range("M136").Value = -23097178.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1706.5625
$ws.Range("I3").Value = 916
$ws.Range("K3").Value = 916
$ws.Range("M3").Value = -802
$ws.Range("H37").Value = 5849.5
$ws.Range("I37").Value = 1699
$ws.Range("K37").Value = 1699
$ws.Range("M37").Value = -1562
$ws.Range("H95").Value = 42704.332
$ws.Range("J95").Value = 42704.332
$ws.Range("L95").Value = 42704.332
$ws.Range("N95").Value = -48196.332
$ws.Range("H105").Value = 1147646.8
$ws.Range("I105").Value = 1527176.5
$ws.Range("K105").Value = 1527176.5
$ws.Range("M105").Value = -1525429.5
$ws.Range("H107").Value = 7351.3335
$ws.Range("I107").Value = 8181.6
$ws.Range("K107").Value = 8181.6
$ws.Range("M107").Value = -6261.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 10.375
$ws.Range("J7").Value = 11.5
$ws.Range("L7").Value = 11.5
$ws.Range("N7").Value = -237.5
$ws.Range("H15").Value = 3699.6667
$ws.Range("J15").Value = 3699.6667
$ws.Range("L15").Value = 3699.6667
$ws.Range("N15").Value = -4039.6667
$ws.Range("H31").Value = 34486864
$ws.Range("I31").Value = 71432790
$ws.Range("K31").Value = 71432790
$ws.Range("M31").Value = -71432495
$ws.Range("H34").Value = 34486864
$ws.Range("I34").Value = 71432790
$ws.Range("K34").Value = 71432790
$ws.Range("M34").Value = -71432588
$ws.Range("H58").Value = 1829.289
$ws.Range("I58").Value = 1979.5
$ws.Range("J58").Value = 1623.7368
$ws.Range("K58").Value = 1979.5
$ws.Range("L58").Value = 1623.7368
$ws.Range("M58").Value = -1776.5
$ws.Range("N58").Value = -2029.7368
$ws.Range("H62").Value = 14499741
$ws.Range("I62").Value = 7035.5713
$ws.Range("J62").Value = 37043948
$ws.Range("K62").Value = 7035.5713
$ws.Range("L62").Value = 37043948
$ws.Range("M62").Value = -6411.5713
$ws.Range("N62").Value = -37045196
$ws.Range("H65").Value = 14499741
$ws.Range("I65").Value = 7035.5713
$ws.Range("J65").Value = 37043948
$ws.Range("K65").Value = 35177.85649999999
$ws.Range("L65").Value = 185219740
$ws.Range("M65").Value = -32057.85649999999
$ws.Range("N65").Value = -185225980
$ws.Range("H99").Value = 12275.571
$ws.Range("I99").Value = 7863.9375
$ws.Range("K99").Value = 7863.9375
$ws.Range("M99").Value = -6365.9375
$ws.Range("H105").Value = 1296.6666
$ws.Range("I105").Value = 1083.75
$ws.Range("K105").Value = 1083.75
$ws.Range("M105").Value = 663.25
$ws.Range("H107").Value = 994.05
$ws.Range("I107").Value = 504.08334
$ws.Range("J107").Value = 1729
$ws.Range("K107").Value = 504.08334
$ws.Range("L107").Value = 1729
$ws.Range("M107").Value = 1415.91666
$ws.Range("N107").Value = -5569
$ws.Range("H122").Value = 3958.2727
$ws.Range("I122").Value = 3567.75
$ws.Range("J122").Value = 4999.6665
$ws.Range("K122").Value = 10703.25
$ws.Range("L122").Value = 14998.9995
$ws.Range("M122").Value = -8253.25
$ws.Range("N122").Value = -19898.9995
$ws.Range("H126").Value = 12275.571
$ws.Range("I126").Value = 7863.9375
$ws.Range("K126").Value = 23591.8125
$ws.Range("M126").Value = -21121.8125
$ws.Range("H136").Value = 1829.289
$ws.Range("I136").Value = 1979.5
$ws.Range("J136").Value = 1623.7368
$ws.Range("K136").Value = 5938.5
$ws.Range("L136").Value = 4871.2104
$ws.Range("M136").Value = -3388.5
$ws.Range("N136").Value = -9971.2104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 380.76923
$ws.Range("I5").Value = 296.77777
$ws.Range("J5").Value = 569.75
$ws.Range("K5").Value = 890.33331
$ws.Range("L5").Value = 1709.25
$ws.Range("M5").Value = -778.33331
$ws.Range("N5").Value = -1933.25
$ws.Range("H33").Value = 6929740
$ws.Range("I33").Value = 95
$ws.Range("J33").Value = 8662151
$ws.Range("K33").Value = 570
$ws.Range("L33").Value = 51972906
$ws.Range("M33").Value = -287
$ws.Range("N33").Value = -51973472
$ws.Range("H38").Value = 21.777779
$ws.Range("I38").Value = 14.166667
$ws.Range("K38").Value = 42.500001
$ws.Range("M38").Value = 304.499999
$ws.Range("H62").Value = 16036
$ws.Range("J62").Value = 16036
$ws.Range("L62").Value = 48108
$ws.Range("N62").Value = -49480
$ws.Range("H63").Value = 20048.47
$ws.Range("I63").Value = 13100.6
$ws.Range("J63").Value = 22943.416
$ws.Range("K63").Value = 39301.8
$ws.Range("L63").Value = 68830.24800000001
$ws.Range("M63").Value = -38552.8
$ws.Range("N63").Value = -70328.24800000001
$ws.Range("H65").Value = 16036
$ws.Range("J65").Value = 16036
$ws.Range("L65").Value = 144324
$ws.Range("N65").Value = -151188
$ws.Range("H66").Value = 20048.47
$ws.Range("I66").Value = 13100.6
$ws.Range("J66").Value = 22943.416
$ws.Range("K66").Value = 117905.4
$ws.Range("L66").Value = 206490.744
$ws.Range("M66").Value = -114161.4
$ws.Range("N66").Value = -213978.744
$ws.Range("H69").Value = 12193
$ws.Range("J69").Value = 14424.333
$ws.Range("L69").Value = 43272.999
$ws.Range("N69").Value = -44894.999
$ws.Range("H72").Value = 12193
$ws.Range("J72").Value = 14424.333
$ws.Range("L72").Value = 129818.997
$ws.Range("N72").Value = -137930.997
$ws.Range("H97").Value = 2161.6667
$ws.Range("I97").Value = 2332.6667
$ws.Range("J97").Value = 2104.6667
$ws.Range("K97").Value = 6998.000100000001
$ws.Range("L97").Value = 6314.000100000001
$ws.Range("M97").Value = -6502.000100000001
$ws.Range("N97").Value = -7306.000100000001
$ws.Range("H104").Value = 33333
$ws.Range("J104").Value = 33333
$ws.Range("L104").Value = 99999
$ws.Range("N104").Value = -105241
$ws.Range("H107").Value = 5693011.5
$ws.Range("J107").Value = 9107328
$ws.Range("L107").Value = 27321984
$ws.Range("N107").Value = -27325824
$ws.Range("H131").Value = 4439.864
$ws.Range("J131").Value = 10315.5
$ws.Range("L131").Value = 30946.5
$ws.Range("N131").Value = -41026.5
$ws.Range("H135").Value = 380.76923
$ws.Range("I135").Value = 296.77777
$ws.Range("J135").Value = 569.75
$ws.Range("K135").Value = 2670.99993
$ws.Range("L135").Value = 5127.75
$ws.Range("M135").Value = -135.9999299999999
$ws.Range("N135").Value = -10197.75
$ws.Range("H140").Value = 6972.3
$ws.Range("I140").Value = 4890.9165
$ws.Range("K140").Value = 14672.7495
$ws.Range("M140").Value = -9492.749500000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 67500
$ws.Range("J93").Value = 67500
$ws.Range("L93").Value = 67500
$ws.Range("N93").Value = -71244
$ws.Range("H98").Value = 23154.166
$ws.Range("J98").Value = 23154.166
$ws.Range("L98").Value = 23154.166
$ws.Range("N98").Value = -29144.166
$ws.Range("H113").Value = 1687436.5
$ws.Range("J113").Value = 3706902.2
$ws.Range("L113").Value = 3706902.2
$ws.Range("N113").Value = -3711242.2
$ws.Range("H122").Value = 3420.9524
$ws.Range("I122").Value = 3386.4614
$ws.Range("J122").Value = 3477
$ws.Range("K122").Value = 10159.3842
$ws.Range("L122").Value = 10431
$ws.Range("M122").Value = -7709.3842
$ws.Range("N122").Value = -15331
$ws.Range("H126").Value = 3278.8333
$ws.Range("I126").Value = 2918.25
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 8754.75
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -6284.75
$ws.Range("N126").Value = -16940
$ws.Range("H132").Value = 6161452
$ws.Range("I132").Value = 3227.8462
$ws.Range("K132").Value = 9683.5386
$ws.Range("M132").Value = -7153.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9450.895
$ws.Range("I7").Value = 8120.615
$ws.Range("K7").Value = 8120.615
$ws.Range("M7").Value = -8008.615
$ws.Range("H55").Value = 989.0417
$ws.Range("I55").Value = 401.8
$ws.Range("J55").Value = 1408.5
$ws.Range("K55").Value = 401.8
$ws.Range("L55").Value = 1408.5
$ws.Range("M55").Value = -228.8
$ws.Range("N55").Value = -1754.5
$ws.Range("H68").Value = 8334977
$ws.Range("I68").Value = 10418166
$ws.Range("K68").Value = 10418166
$ws.Range("M68").Value = -10417417
$ws.Range("H71").Value = 8334977
$ws.Range("I71").Value = 10418166
$ws.Range("K71").Value = 52090830
$ws.Range("M71").Value = -52087086
$ws.Range("H126").Value = 9450.895
$ws.Range("I126").Value = 8120.615
$ws.Range("K126").Value = 24361.845
$ws.Range("M126").Value = -21891.845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 16779.166
$ws.Range("I74").Value = 33600
$ws.Range("K74").Value = 33600
$ws.Range("M74").Value = -32664
$ws.Range("H77").Value = 16779.166
$ws.Range("I77").Value = 33600
$ws.Range("K77").Value = 100800
$ws.Range("M77").Value = -96120
$ws.Range("H100").Value = 850.95654
$ws.Range("I100").Value = 881.25
$ws.Range("J100").Value = 781.7143
$ws.Range("K100").Value = 1762.5
$ws.Range("L100").Value = 1563.4286
$ws.Range("M100").Value = -1221.5
$ws.Range("N100").Value = -2645.4286
$ws.Range("H122").Value = 2809.4211
$ws.Range("I122").Value = 2817.6
$ws.Range("K122").Value = 8452.799999999999
$ws.Range("M122").Value = -6002.799999999999
$ws.Range("H126").Value = 7463.636
$ws.Range("I126").Value = 7442.5264
$ws.Range("J126").Value = 7597.3335
$ws.Range("K126").Value = 22327.5792
$ws.Range("L126").Value = 22792.0005
$ws.Range("M126").Value = -19857.5792
$ws.Range("N126").Value = -27732.0005
$ws.Range("H136").Value = 264003.25
$ws.Range("I136").Value = 832.3333
$ws.Range("J136").Value = 2000931.4
$ws.Range("K136").Value = 2496.9999
$ws.Range("L136").Value = 6002794.199999999
$ws.Range("M136").Value = 53.0001000000002
$ws.Range("N136").Value = -6007894.199999999
